$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-19 Monday" "2024-08-20 Tuesday"

Replace-Text "97×50=4850" "94×11=1034"
Replace-Text "90×16=1440" "26×94=2444"
Replace-Text "36×13=468" "59×63=3717"
Replace-Text "82×64=5248" "22×49=1078"
Replace-Text "88×47=4136" "91×71=6461"

Replace-Text "33×49=1617" "42×54=2268"
Replace-Text "46×52=2392" "36×31=1116"
Replace-Text "72×43=3096" "78×44=3432"
Replace-Text "49×46=2254" "54×62=3348"
Replace-Text "76×74=5624" "46×14=644"

Replace-Text "25×68=1700" "75×67=5025"
Replace-Text "22×44=968" "32×31=992"
Replace-Text "44×36=1584" "93×93=8649"
Replace-Text "15×22=330" "69×99=6831"
Replace-Text "58×66=3828" "84×71=5964"

Replace-Text "48×97=4656" "49×64=3136"
Replace-Text "98×46=4508" "27×55=1485"
Replace-Text "15×13=195" "13×99=1287"
Replace-Text "13×45=585" "85×45=3825"
Replace-Text "30×27=810" "54×53=2862"

Replace-Text "11×77=847" "25×75=1875"
Replace-Text "56×38=2128" "98×52=5096"
Replace-Text "85×59=5015" "30×82=2460"
Replace-Text "21×75=1575" "25×50=1250"
Replace-Text "99×70=6930" "96×19=1824"
